# Update the sample data values on Sheet1 to match the regression-test
# fixture's new expected numbers. The cells already carry the correct
# styles (B column -> currency style, C column -> thousands-style), so we
# only need to change the stored values, not the formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 141.5
$ws.Range("C2").Value = 92

$ws.Range("B3").Value = 0.314
$ws.Range("C3").Value = 15

$ws.Range("B4").Value = 653.5
$ws.Range("C4").Value = 14
